$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: theta_se standard errors
$ws.Range("B4").Value = "(0.79)"
$ws.Range("C4").Value = "(0.7)"
$ws.Range("D4").Value = "(0.92)"
$ws.Range("E4").Value = "(0.89)"
$ws.Range("F4").Value = "(0.96)"
$ws.Range("G4").Value = "(0.9)"
$ws.Range("H4").Value = "(1.04)"
$ws.Range("I4").Value = "(0.94)"
$ws.Range("J4").Value = "(0.95)"
$ws.Range("K4").Value = "(0.65)"
$ws.Range("L4").Value = "(1.11)"

# Row 6: lambda_se standard errors
$ws.Range("B6").Value = "(0.5)"
$ws.Range("C6").Value = "(0.46)"
$ws.Range("D6").Value = "(0.59)"
$ws.Range("E6").Value = "(0.68)"
$ws.Range("F6").Value = "(0.64)"
$ws.Range("G6").Value = "(0.64)"
$ws.Range("H6").Value = "(0.8)"
$ws.Range("I6").Value = "(0.73)"
$ws.Range("J6").Value = "(0.67)"
$ws.Range("K6").Value = "(0.56)"
$ws.Range("L6").Value = "(0.86)"
